$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("Período") before the existing "Categoria" column
$ws.Columns.Item(2).Insert()

# Update header row
$ws.Range("B1").Value = "Período"
$ws.Range("C1").Value = "Categoria"
$ws.Range("D1").Value = "Valor"

# Data for the "Período" (new column B), "Categoria" (column C) and "Valor" (column D)
$periodos = @(
  "Quantidade 2024/2014","Quantidade 2024/2014","Quantidade 2024/2014","Quantidade 2024/2014","Quantidade 2024/2014","Quantidade 2024/2014",
  "Valor 2024/2014","Valor 2024/2014","Valor 2024/2014","Valor 2024/2014","Valor 2024/2014","Valor 2024/2014",
  "Quantidade 2024/2023","Quantidade 2024/2023","Quantidade 2024/2023","Quantidade 2024/2023","Quantidade 2024/2023","Quantidade 2024/2023",
  "Valor 2024/2023","Valor 2024/2023","Valor 2024/2023","Valor 2024/2023","Valor 2024/2023","Valor 2024/2023"
)

$categorias = @(
  "Variação em dez anos","Variação em dez anos","Variação em dez anos","Variação em dez anos","Variação em dez anos","Variação em dez anos",
  "Variação em dez anos","Variação em dez anos","Variação em dez anos","Variação em dez anos","Variação em dez anos","Variação em dez anos",
  "Variação do último ano","Variação do último ano","Variação do último ano","Variação do último ano","Variação do último ano","Variação do último ano",
  "Variação do último ano","Variação do último ano","Variação do último ano","Variação do último ano","Variação do último ano","Variação do último ano"
)

$valores = @(
  67.40443935635693, 68.53318765364655, -9.350826619328432, -59.33206703373326, 637.4722838137473, 212.9210042865891,
  119.0594447084194, 42.53980162572068, -20.88500066046327, -52.33362673955766, 661.7192065911167, 76.15988489401225,
  -18.40979541227526, 1.224477620844457, 11.20181579898588, 1.829587665484669, 34.65587044534413, -6.41025641025641,
  2.934781326149883, -36.31319223224396, 21.18909609656381, 4.073232865106783, 38.51361536014834, -7.859334322257486
)

for ($i = 0; $i -lt 24; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 2).Value = $periodos[$i]
  $ws.Cells.Item($row, 3).Value = $categorias[$i]
  $ws.Cells.Item($row, 4).Value = $valores[$i]
}
